# Update "想去人数" (F column) figures across the sheets, as scraped at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value  = 857
$ws1.Cells.Item(3, 6).Value  = 13768
$ws1.Cells.Item(4, 6).Value  = 13558
$ws1.Cells.Item(5, 6).Value  = 1050
$ws1.Cells.Item(6, 6).Value  = 807
$ws1.Cells.Item(11, 6).Value = 53
$ws1.Cells.Item(12, 6).Value = 760
$ws1.Cells.Item(13, 6).Value = 2142
$ws1.Cells.Item(14, 6).Value = 96
$ws1.Cells.Item(17, 6).Value = 119
$ws1.Cells.Item(19, 6).Value = 525
$ws1.Cells.Item(20, 6).Value = 431
$ws1.Cells.Item(21, 6).Value = 392
$ws1.Cells.Item(22, 6).Value = 321
$ws1.Cells.Item(23, 6).Value = 261
$ws1.Cells.Item(24, 6).Value = 831
$ws1.Cells.Item(25, 6).Value = 81

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(6, 6).Value  = 165
$ws2.Cells.Item(7, 6).Value  = 1491
$ws2.Cells.Item(10, 6).Value = 23
$ws2.Cells.Item(11, 6).Value = 66

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 222
$ws3.Cells.Item(3, 6).Value = 106

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value  = 222
$ws4.Cells.Item(3, 6).Value  = 857
$ws4.Cells.Item(4, 6).Value  = 13768
$ws4.Cells.Item(5, 6).Value  = 13558
$ws4.Cells.Item(6, 6).Value  = 1050
$ws4.Cells.Item(7, 6).Value  = 807
$ws4.Cells.Item(12, 6).Value = 53
$ws4.Cells.Item(13, 6).Value = 760
$ws4.Cells.Item(16, 6).Value = 2142
$ws4.Cells.Item(17, 6).Value = 96
$ws4.Cells.Item(20, 6).Value = 119
$ws4.Cells.Item(24, 6).Value = 106
$ws4.Cells.Item(25, 6).Value = 106
$ws4.Cells.Item(26, 6).Value = 525
$ws4.Cells.Item(27, 6).Value = 431
$ws4.Cells.Item(28, 6).Value = 393
$ws4.Cells.Item(29, 6).Value = 321
$ws4.Cells.Item(30, 6).Value = 261
$ws4.Cells.Item(31, 6).Value = 831
$ws4.Cells.Item(32, 6).Value = 165
$ws4.Cells.Item(33, 6).Value = 1491
$ws4.Cells.Item(36, 6).Value = 23
$ws4.Cells.Item(37, 6).Value = 81
$ws4.Cells.Item(38, 6).Value = 66
